$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 10:05"

# --- Refreshed COVID-19 country statistics ---
# A handful of countries swap rank (and so swap rows) as updated case
# counts push them past/behind a neighbour in the table; those rows get
# a new country name together with new figures, while the row that used
# to hold that country now shows the next one down (with ITS figures).

# Row 6
$ws.Range("B6").Value = 242271
$ws.Range("C6").Value = 10028
$ws.Range("D6").Value = 48003
$ws.Range("E6").Value = 192056
$ws.Range("G6").Value = 96
$ws.Range("H6").Value = 2212

# Row 25
$ws.Range("F25").Value = 79

# Row 51
$ws.Range("B51").Value = 8223
$ws.Range("C51").Value = 25
$ws.Range("D51").Value = 4899
$ws.Range("E51").Value = 3040
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 284

# Row 65
$ws.Range("A65").Value = "Oman"
$ws.Range("B65").Value = 4019
$ws.Range("C65").Value = 298
$ws.Range("D65").Value = 1289
$ws.Range("E65").Value = 2713
$ws.Range("F65").Value = 27
$ws.Range("H65").Value = 17

# Row 66
$ws.Range("A66").Value = "Luxemburgo"
$ws.Range("B66").Value = 3894
$ws.Range("D66").Value = 3610
$ws.Range("E66").Value = 182
$ws.Range("F66").Value = 22
$ws.Range("H66").Value = 102

# Row 73
$ws.Range("D73").Value = 1735
$ws.Range("E73").Value = 925

# Row 88
$ws.Range("A88").Value = "Lituania"
$ws.Range("B88").Value = 1505
$ws.Range("C88").Value = 14
$ws.Range("D88").Value = 908
$ws.Range("E88").Value = 543
$ws.Range("F88").Value = 17
$ws.Range("G88").Value = 4
$ws.Range("H88").Value = 54

# Row 89
$ws.Range("A89").Value = "Nueva Zelanda"
$ws.Range("B89").Value = 1497
$ws.Range("D89").Value = 1402
$ws.Range("E89").Value = 74
$ws.Range("F89").Value = 2
$ws.Range("H89").Value = 21

# Row 90
$ws.Range("B90").Value = 1469
$ws.Range("C90").Value = 4
$ws.Range("D90").Value = 1060
$ws.Range("E90").Value = 382

# Row 193
$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# Row 194
$ws.Range("A194").Value = "Nueva Caledonia"
$ws.Range("D194").Value = 18
$ws.Range("H194").Value = 0

# Row 198
$ws.Range("A198").Value = "Curazao"
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 1

# Row 199
$ws.Range("A199").Value = "Dominica"
$ws.Range("D199").Value = 15
$ws.Range("H199").Value = 0

# Row 214
$ws.Range("A214").Value = "Sahara Occidental"
$ws.Range("D214").Value = 6
$ws.Range("E214").Value = 0

# Row 216
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
